$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 454634.47
$ws.Range("I6").Value = 500092.9
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 1500278.7
$ws.Range("L6").Value = 150
$ws.Range("M6").Value = -1500166.7
$ws.Range("N6").Value = -374

$ws.Range("H28").Value = 949.5
$ws.Range("I28").Value = 899
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 899
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -414
$ws.Range("N28").Value = -1970

$ws.Range("H31").Value = 243.71428
$ws.Range("I31").Value = 243.71428
$ws.Range("K31").Value = 731.14284
$ws.Range("M31").Value = -501.14284

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2925.8
$ws.Range("I32").Value = 2925.8
$ws.Range("K32").Value = 2925.8
$ws.Range("M32").Value = -2638.8

$ws.Range("H97").Value = 532.5714
$ws.Range("I97").Value = 532.5714
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 532.5714
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -36.57140000000004
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 5166.3335
$ws.Range("I102").Value = 5166.3335
$ws.Range("K102").Value = 5166.3335
$ws.Range("M102").Value = -3544.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1134.1666
$ws.Range("I86").Value = 899.5
$ws.Range("J86").Value = 1603.5
$ws.Range("K86").Value = 899.5
$ws.Range("L86").Value = 1603.5
$ws.Range("M86").Value = 223.5
$ws.Range("N86").Value = -3849.5

$ws.Range("H89").Value = 1134.1666
$ws.Range("I89").Value = 899.5
$ws.Range("J89").Value = 1603.5
$ws.Range("K89").Value = 4497.5
$ws.Range("L89").Value = 8017.5
$ws.Range("M89").Value = 1118.5
$ws.Range("N89").Value = -19249.5

$ws.Range("H105").Value = 1699.75
$ws.Range("I105").Value = 1599.6666
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1599.6666
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 147.3334
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 201.2683
$ws.Range("I7").Value = 244.11111
$ws.Range("J7").Value = 189.21875
$ws.Range("K7").Value = 244.11111
$ws.Range("L7").Value = 189.21875
$ws.Range("M7").Value = -131.11111
$ws.Range("N7").Value = -415.21875

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H55").Value = 12995
$ws.Range("J55").Value = 12995
$ws.Range("L55").Value = 12995
$ws.Range("N55").Value = -13625

$ws.Range("H92").Value = 34999.5
$ws.Range("J92").Value = 34999.5
$ws.Range("L92").Value = 34999.5
$ws.Range("N92").Value = -39991.5

$ws.Range("H105").Value = 2398.5
$ws.Range("I105").Value = 2398.5
$ws.Range("K105").Value = 2398.5
$ws.Range("M105").Value = -651.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 805
$ws.Range("I2").Value = 314.125
$ws.Range("K2").Value = 1884.75
$ws.Range("M2").Value = -1771.75

$ws.Range("H17").Value = 47.333332
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 46
$ws.Range("K17").Value = 150
$ws.Range("L17").Value = 138
$ws.Range("M17").Value = 19
$ws.Range("N17").Value = -476

$ws.Range("H68").Value = 999
$ws.Range("I68").Value = 999
$ws.Range("K68").Value = 2997
$ws.Range("M68").Value = -2186

$ws.Range("H71").Value = 999
$ws.Range("I71").Value = 999
$ws.Range("K71").Value = 8991
$ws.Range("M71").Value = -4935

$ws.Range("H98").Value = 1126.5
$ws.Range("I98").Value = 1169
$ws.Range("K98").Value = 3507
$ws.Range("M98").Value = -2009

$ws.Range("H104").Value = 756
$ws.Range("I104").Value = 1013
$ws.Range("J104").Value = 499
$ws.Range("K104").Value = 3039
$ws.Range("L104").Value = 1497
$ws.Range("M104").Value = -418
$ws.Range("N104").Value = -6739

$ws.Range("H131").Value = 6070.1665
$ws.Range("I131").Value = 6070.1665
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 18210.4995
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -13170.4995
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 45.2
$ws.Range("I2").Value = 25.25
$ws.Range("K2").Value = 25.25
$ws.Range("M2").Value = 87.75

$ws.Range("H13").Value = 109.25
$ws.Range("I13").Value = 109.25
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 109.25
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 29.75
$ws.Range("N13").ClearContents()

$ws.Range("H92").Value = 2928.5
$ws.Range("J92").Value = 2928.5
$ws.Range("L92").Value = 2928.5
$ws.Range("N92").Value = -6672.5

$ws.Range("H107").Value = 3789.25
$ws.Range("I107").Value = 5748.5
$ws.Range("J107").Value = 1830
$ws.Range("K107").Value = 5748.5
$ws.Range("L107").Value = 1830
$ws.Range("M107").Value = -3828.5
$ws.Range("N107").Value = -5670

$ws.Range("H113").Value = 311
$ws.Range("I113").Value = 311
$ws.Range("K113").Value = 311
$ws.Range("M113").Value = 1859

$ws.Range("H122").Value = 5191.8
$ws.Range("I122").Value = 1490.25
$ws.Range("J122").Value = 19998
$ws.Range("K122").Value = 4470.75
$ws.Range("L122").Value = 59994
$ws.Range("M122").Value = -2020.75
$ws.Range("N122").Value = -64894

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2169.3333
$ws.Range("I16").Value = 2088
$ws.Range("J16").Value = 2210
$ws.Range("K16").Value = 2088
$ws.Range("L16").Value = 2210
$ws.Range("M16").Value = -1918
$ws.Range("N16").Value = -2550

$ws.Range("H100").Value = 1519
$ws.Range("I100").Value = 1519
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1519
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -978
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H132").Value = 4884.375
$ws.Range("I132").Value = 3346
$ws.Range("K132").Value = 10038
$ws.Range("M132").Value = -7508
